$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.898.01"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.231.40"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'249.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.74%  "
$ws.Range("D6").Value = "'0.630"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").Value = "'71.73"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.76%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.598"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.13%  "
$ws.Range("D10").Value = "'40.65"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +15.13%  "
$ws.Range("D11").Value = "'0.0972"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").Value = "'58.05"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'7.15"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.36%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "2.567.84"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "'14.98"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'0.865"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "2.225.42"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "41.959.35"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "0.0₃0971"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'6.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "'73.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'235.74"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.32%  "
$ws.Range("D25").Value = "'4.03"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.91%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'2.53"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.38%  "
$ws.Range("D28").Value = "'10.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.13%  "
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "'171.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "'20.75"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "'0.124"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("D33").Value = "'0.125"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").Value = "'4.74"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").Value = "'26.53"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +25.64%  "
$ws.Range("D38").Value = "'4.01"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +11.18%  "
$ws.Range("D39").Value = "'0.0297"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +11.16%  "
$ws.Range("D40").Value = "'2.29"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("D41").Value = "'6.02"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").Value = "'66.54"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "'12.09"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +22.62%  "
$ws.Range("D44").Value = "'0.206"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.71%  "
$ws.Range("D45").Value = "'4.96"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("D46").Value = "'8.79"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "'4.66"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.47%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'1.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.32%  "
$ws.Range("D51").Value = "'1.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.15%  "
